# PanelApp "Cerebral malformations" sheet update
# - Insert a new gene row (VPS50) after VLDLR, pushing AKT2..SLC35A2 down by one row
# - Add a new "time_taken" metadata column (F) with a timestamp for every data row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Create room for the new 153rd row by duplicating the formatting that
#    the static row-index column (A) already uses, then give it its value.
# ---------------------------------------------------------------------------
$ws.Range("A152").Copy()
$ws.Range("A153").PasteSpecial(-4122)
$ws.Range("A153").Value = 151

# ---------------------------------------------------------------------------
# 2. Re-write columns B (geneSymbol), C (geneName) and D (geneConfidence) for
#    rows 124-153 so that the VPS50 entry is inserted right after VLDLR and
#    every subsequent gene moves down by one row. Column A is left untouched
#    (it is just the static 0-based row index already present in the sheet).
# ---------------------------------------------------------------------------
$shiftedRows = @(
  @("VPS50", "VPS50, EARP/GARPII complex subunit", "2"),
  @("AKT2", "AKT serine/threonine kinase 2", "1"),
  @("EOMES", "eomesodermin", "1"),
  @("ERMARD", "ER membrane associated RNA degradation", "1"),
  @("FOXH1", "forkhead box H1", "1"),
  @("GCM2", "glial cells missing homolog 2", "1"),
  @("GMPPB", "GDP-mannose pyrophosphorylase B", "1"),
  @("HRAS", "HRas proto-oncogene, GTPase", "1"),
  @("KRAS", "KRAS proto-oncogene, GTPase", "1"),
  @("MCF2", "MCF.2 cell line derived transforming sequence", "1"),
  @("MTOR", "mechanistic target of rapamycin kinase", "1"),
  @("NODAL", "nodal growth differentiation factor", "1"),
  @("NRAS", "NRAS proto-oncogene, GTPase", "1"),
  @("PEX11A", "peroxisomal biogenesis factor 11 alpha", "1"),
  @("POMK", "protein-O-mannose kinase", "1"),
  @("SMAD2", "SMAD family member 2", "1"),
  @("TBC1D7", "TBC1 domain family member 7", "1"),
  @("TSC2", "TSC complex subunit 2", "1"),
  @("ATP1A3", "ATPase Na+/K+ transporting subunit alpha 3", "0"),
  @("ENO1", "enolase 1", "0"),
  @("GRIN2B", "glutamate ionotropic receptor NMDA type subunit 2B", "0"),
  @("MAPK8IP3", "mitogen-activated protein kinase 8 interacting protein 3", "0"),
  @("NPRL2", "NPR2 like, GATOR1 complex subunit", "0"),
  @("NPRL3", "NPR3 like, GATOR1 complex subunit", "0"),
  @("PTEN", "phosphatase and tensin homolog", "0"),
  @("RAB18", "RAB18, member RAS oncogene family", "0"),
  @("RAB3GAP1", "RAB3 GTPase activating protein catalytic subunit 1", "0"),
  @("RAB3GAP2", "RAB3 GTPase activating non-catalytic protein subunit 2", "0"),
  @("SCN3A", "sodium voltage-gated channel alpha subunit 3", "0"),
  @("SLC35A2", "solute carrier family 35 member A2", "0")
)

$bcRange = $ws.Range("B124:C153")
$bcArr = New-Object 'object[,]' 30,2
for ($i = 0; $i -lt 30; $i++) {
    $bcArr[$i,0] = $shiftedRows[$i][0]
    $bcArr[$i,1] = $shiftedRows[$i][1]
}
$bcRange.Value = $bcArr

# Column D holds numeric-looking confidence codes that must stay text, as in
# the rest of the column.
$dRange = $ws.Range("D124:D153")
$dRange.NumberFormat = "@"
$dArr = New-Object 'object[,]' 30,1
for ($i = 0; $i -lt 30; $i++) {
    $dArr[$i,0] = $shiftedRows[$i][2]
}
$dRange.Value = $dArr
# Reset the visual style back to the plain/default one used elsewhere in
# column D, while keeping the cached value stored as text.
$ws.Range("E124").Copy()
$dRange.PasteSpecial(-4122)

# Column E (panel) for the brand new row 153.
$ws.Range("E153").Value = "Cerebral malformations"

# ---------------------------------------------------------------------------
# 3. Add the new "time_taken" column (F) with header + one timestamp per row.
# ---------------------------------------------------------------------------
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

$f_values = @(
  "2021-10-05 13:38:58.912983",
  "2021-10-05 13:38:58.912995",
  "2021-10-05 13:38:58.912998",
  "2021-10-05 13:38:58.913001",
  "2021-10-05 13:38:58.913004",
  "2021-10-05 13:38:58.913007",
  "2021-10-05 13:38:58.913009",
  "2021-10-05 13:38:58.913011",
  "2021-10-05 13:38:58.913014",
  "2021-10-05 13:38:58.913017",
  "2021-10-05 13:38:58.913019",
  "2021-10-05 13:38:58.913022",
  "2021-10-05 13:38:58.913024",
  "2021-10-05 13:38:58.913027",
  "2021-10-05 13:38:58.913030",
  "2021-10-05 13:38:58.913032",
  "2021-10-05 13:38:58.913035",
  "2021-10-05 13:38:58.913038",
  "2021-10-05 13:38:58.913040",
  "2021-10-05 13:38:58.913043",
  "2021-10-05 13:38:58.913045",
  "2021-10-05 13:38:58.913048",
  "2021-10-05 13:38:58.913050",
  "2021-10-05 13:38:58.913053",
  "2021-10-05 13:38:58.913056",
  "2021-10-05 13:38:58.913058",
  "2021-10-05 13:38:58.913061",
  "2021-10-05 13:38:58.913063",
  "2021-10-05 13:38:58.913066",
  "2021-10-05 13:38:58.913068",
  "2021-10-05 13:38:58.913071",
  "2021-10-05 13:38:58.913073",
  "2021-10-05 13:38:58.913076",
  "2021-10-05 13:38:58.913078",
  "2021-10-05 13:38:58.913081",
  "2021-10-05 13:38:58.913083",
  "2021-10-05 13:38:58.913086",
  "2021-10-05 13:38:58.913089",
  "2021-10-05 13:38:58.913091",
  "2021-10-05 13:38:58.913094",
  "2021-10-05 13:38:58.913097",
  "2021-10-05 13:38:58.913099",
  "2021-10-05 13:38:58.913102",
  "2021-10-05 13:38:58.913104",
  "2021-10-05 13:38:58.913107",
  "2021-10-05 13:38:58.913109",
  "2021-10-05 13:38:58.913112",
  "2021-10-05 13:38:58.913114",
  "2021-10-05 13:38:58.913117",
  "2021-10-05 13:38:58.913119",
  "2021-10-05 13:38:58.913122",
  "2021-10-05 13:38:58.913125",
  "2021-10-05 13:38:58.913127",
  "2021-10-05 13:38:58.913130",
  "2021-10-05 13:38:58.913133",
  "2021-10-05 13:38:58.913135",
  "2021-10-05 13:38:58.913138",
  "2021-10-05 13:38:58.913140",
  "2021-10-05 13:38:58.913143",
  "2021-10-05 13:38:58.913145",
  "2021-10-05 13:38:58.913148",
  "2021-10-05 13:38:58.913150",
  "2021-10-05 13:38:58.913153",
  "2021-10-05 13:38:58.913155",
  "2021-10-05 13:38:58.913159",
  "2021-10-05 13:38:58.913162",
  "2021-10-05 13:38:58.913164",
  "2021-10-05 13:38:58.913167",
  "2021-10-05 13:38:58.913169",
  "2021-10-05 13:38:58.913172",
  "2021-10-05 13:38:58.913174",
  "2021-10-05 13:38:58.913177",
  "2021-10-05 13:38:58.913179",
  "2021-10-05 13:38:58.913182",
  "2021-10-05 13:38:58.913184",
  "2021-10-05 13:38:58.913187",
  "2021-10-05 13:38:58.913191",
  "2021-10-05 13:38:58.913194",
  "2021-10-05 13:38:58.913197",
  "2021-10-05 13:38:58.913200",
  "2021-10-05 13:38:58.913202",
  "2021-10-05 13:38:58.913204",
  "2021-10-05 13:38:58.913207",
  "2021-10-05 13:38:58.913209",
  "2021-10-05 13:38:58.913212",
  "2021-10-05 13:38:58.913214",
  "2021-10-05 13:38:58.913217",
  "2021-10-05 13:38:58.913219",
  "2021-10-05 13:38:58.913222",
  "2021-10-05 13:38:58.913224",
  "2021-10-05 13:38:58.913227",
  "2021-10-05 13:38:58.913230",
  "2021-10-05 13:38:58.913233",
  "2021-10-05 13:38:58.913236",
  "2021-10-05 13:38:58.913239",
  "2021-10-05 13:38:58.913241",
  "2021-10-05 13:38:58.913244",
  "2021-10-05 13:38:58.913246",
  "2021-10-05 13:38:58.913249",
  "2021-10-05 13:38:58.913251",
  "2021-10-05 13:38:58.913254",
  "2021-10-05 13:38:58.913256",
  "2021-10-05 13:38:58.913259",
  "2021-10-05 13:38:58.913261",
  "2021-10-05 13:38:58.913264",
  "2021-10-05 13:38:58.913266",
  "2021-10-05 13:38:58.913269",
  "2021-10-05 13:38:58.913271",
  "2021-10-05 13:38:58.913275",
  "2021-10-05 13:38:58.913278",
  "2021-10-05 13:38:58.913281",
  "2021-10-05 13:38:58.913283",
  "2021-10-05 13:38:58.913286",
  "2021-10-05 13:38:58.913288",
  "2021-10-05 13:38:58.913290",
  "2021-10-05 13:38:58.913293",
  "2021-10-05 13:38:58.913295",
  "2021-10-05 13:38:58.913298",
  "2021-10-05 13:38:58.913300",
  "2021-10-05 13:38:58.913303",
  "2021-10-05 13:38:58.913305",
  "2021-10-05 13:38:58.913308",
  "2021-10-05 13:38:58.913310",
  "2021-10-05 13:38:58.913313",
  "2021-10-05 13:38:58.913315",
  "2021-10-05 13:38:58.913318",
  "2021-10-05 13:38:58.913320",
  "2021-10-05 13:38:58.913323",
  "2021-10-05 13:38:58.913327",
  "2021-10-05 13:38:58.913330",
  "2021-10-05 13:38:58.913333",
  "2021-10-05 13:38:58.913335",
  "2021-10-05 13:38:58.913337",
  "2021-10-05 13:38:58.913340",
  "2021-10-05 13:38:58.913343",
  "2021-10-05 13:38:58.913345",
  "2021-10-05 13:38:58.913348",
  "2021-10-05 13:38:58.913350",
  "2021-10-05 13:38:58.913353",
  "2021-10-05 13:38:58.913355",
  "2021-10-05 13:38:58.913358",
  "2021-10-05 13:38:58.913360",
  "2021-10-05 13:38:58.913363",
  "2021-10-05 13:38:58.913365",
  "2021-10-05 13:38:58.913368",
  "2021-10-05 13:38:58.913370",
  "2021-10-05 13:38:58.913373",
  "2021-10-05 13:38:58.913375",
  "2021-10-05 13:38:58.913378",
  "2021-10-05 13:38:58.913381",
  "2021-10-05 13:38:58.913383",
  "2021-10-05 13:38:58.913386"
)

$fRange = $ws.Range("F2:F153")
$fArr = New-Object 'object[,]' 152,1
for ($i = 0; $i -lt 152; $i++) {
    $fArr[$i,0] = $f_values[$i]
}
$fRange.Value = $fArr

$excel.CutCopyMode = 0

Write-Host "Edit complete"
